$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2170212765957447
$ws.Range("C2").Value = 0.5361702127659574
$ws.Range("J2").Value = 0.00425531914893617
$ws.Range("P2").Value = 0.1319148936170213
$ws.Range("S2").Value = 0.1106382978723404
$ws.Range("C3").Value = 0.04651162790697674
$ws.Range("J3").Value = 0.02325581395348837
$ws.Range("P3").Value = 0.8062015503875969
$ws.Range("S3").Value = 0.124031007751938
$ws.Range("P4").Value = 0.7843137254901961
$ws.Range("S4").Value = 0.2156862745098039
$ws.Range("B6").Value = 0.07589285714285714
$ws.Range("D6").Value = 0.008928571428571428
$ws.Range("F6").Value = 0.09375
$ws.Range("J6").Value = 0.2008928571428572
$ws.Range("O6").Value = 0.02232142857142857
$ws.Range("Q6").Value = 0.15625
$ws.Range("R6").Value = 0.04464285714285714
$ws.Range("S6").Value = 0.3973214285714285
$ws.Range("B7").Value = 0.08629441624365482
$ws.Range("D7").Value = 0.02538071065989848
$ws.Range("F7").Value = 0.04568527918781726
$ws.Range("J7").Value = 0.1370558375634518
$ws.Range("O7").Value = 0.005076142131979695
$ws.Range("Q7").Value = 0.2081218274111675
$ws.Range("R7").Value = 0.06598984771573604
$ws.Range("S7").Value = 0.4263959390862944
$ws.Range("B8").Value = 0.049800796812749
$ws.Range("D8").Value = 0.02191235059760956
$ws.Range("F8").Value = 0.06374501992031872
$ws.Range("J8").Value = 0.1055776892430279
$ws.Range("O8").Value = 0.02191235059760956
$ws.Range("Q8").Value = 0.1673306772908366
$ws.Range("R8").Value = 0.1175298804780877
$ws.Range("S8").Value = 0.4521912350597609
$ws.Range("B9").Value = 0.06722689075630252
$ws.Range("D9").Value = 0.03361344537815126
$ws.Range("F9").Value = 0.04621848739495799
$ws.Range("J9").Value = 0.1050420168067227
$ws.Range("O9").Value = 0.02521008403361345
$ws.Range("Q9").Value = 0.2100840336134454
$ws.Range("R9").Value = 0.1008403361344538
$ws.Range("S9").Value = 0.4117647058823529
$ws.Range("B10").Value = 0.08150943396226415
$ws.Range("D10").Value = 0.01962264150943396
$ws.Range("F10").Value = 0.0649056603773585
$ws.Range("J10").Value = 0.1267924528301887
$ws.Range("O10").Value = 0.01886792452830189
$ws.Range("Q10").Value = 0.2384905660377359
$ws.Range("R10").Value = 0.08452830188679246
$ws.Range("S10").Value = 0.3652830188679245
$ws.Range("G11").Value = 0.1162790697674419
$ws.Range("J11").Value = 0.07364341085271318
$ws.Range("K11").Value = 0.1550387596899225
$ws.Range("L11").Value = 0.6395348837209303
$ws.Range("S11").Value = 0.01550387596899225
$ws.Range("G12").Value = 0.7796610169491526
$ws.Range("J12").Value = 0.1468926553672316
$ws.Range("L12").Value = 0.05649717514124294
$ws.Range("S12").Value = 0.01694915254237288
$ws.Range("F13").Value = 0.02127659574468085
$ws.Range("G13").Value = 0.7659574468085106
$ws.Range("J13").Value = 0.1914893617021277
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.09871244635193133
$ws.Range("I15").Value = 0.06008583690987124
$ws.Range("J15").Value = 0.3819742489270386
$ws.Range("K15").Value = 0.07296137339055794
$ws.Range("M15").Value = 0.01716738197424893
$ws.Range("O15").Value = 0.09442060085836911
$ws.Range("S15").Value = 0.2575107296137339
$ws.Range("F16").Value = 0.01775147928994083
$ws.Range("H16").Value = 0.2366863905325444
$ws.Range("I16").Value = 0.09467455621301775
$ws.Range("J16").Value = 0.3964497041420119
$ws.Range("K16").Value = 0.08284023668639054
$ws.Range("M16").Value = 0.02366863905325444
$ws.Range("O16").Value = 0.04733727810650887
$ws.Range("S16").Value = 0.1005917159763314
$ws.Range("F17").Value = 0.01529636711281071
$ws.Range("H17").Value = 0.2065009560229445
$ws.Range("I17").Value = 0.1051625239005736
$ws.Range("J17").Value = 0.4302103250478012
$ws.Range("K17").Value = 0.08604206500956023
$ws.Range("M17").Value = 0.01338432122370937
$ws.Range("N17").Value = 0.001912045889101338
$ws.Range("O17").Value = 0.05353728489483748
$ws.Range("S17").Value = 0.08795411089866156
$ws.Range("F18").Value = 0.02325581395348837
$ws.Range("H18").Value = 0.1906976744186047
$ws.Range("I18").Value = 0.1581395348837209
$ws.Range("J18").Value = 0.3674418604651163
$ws.Range("K18").Value = 0.09302325581395349
$ws.Range("M18").Value = 0.03255813953488372
$ws.Range("O18").Value = 0.05116279069767442
$ws.Range("S18").Value = 0.08372093023255814
$ws.Range("F19").Value = 0.01427498121712998
$ws.Range("H19").Value = 0.2208865514650639
$ws.Range("I19").Value = 0.0879038317054846
$ws.Range("J19").Value = 0.3786626596543952
$ws.Range("K19").Value = 0.08715251690458302
$ws.Range("M19").Value = 0.01953418482344102
$ws.Range("N19").Value = 0.0007513148009015778
$ws.Range("O19").Value = 0.067618332081142
$ws.Range("S19").Value = 0.1232156273478587
